# Fix NICE LTD SPON ADR price history: correct the open/close/high/low/shares
# and fixed_ticker columns (rows 2-43) which had been mistakenly filled with
# data from other tickers (WIX, NOC, CDNS, ...). Every row should reference NICE.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = 60.43547802289387   # D2 open_price
$ws.Cells.Item(2, 5).Value = 58.68997573852539   # E2 close_price
$ws.Cells.Item(2, 6).Value = 60.51392550424472   # F2 high_price
$ws.Cells.Item(2, 7).Value = 57.02292841702808   # G2 low_price
$ws.Cells.Item(2, 8).Value = 61741703       # H2 shares_outstanding
$ws.Cells.Item(2, 9).Value = "NICE"        # I2 fixed_ticker

$ws.Cells.Item(3, 4).Value = 63.30474971099117   # D3 open_price
$ws.Cells.Item(3, 5).Value = 63.4620246887207   # E3 close_price
$ws.Cells.Item(3, 6).Value = 64.22875801680028   # F3 high_price
$ws.Cells.Item(3, 7).Value = 60.71947711142811   # G3 low_price
$ws.Cells.Item(3, 8).Value = 61741703       # H3 shares_outstanding
$ws.Cells.Item(3, 9).Value = "NICE"        # I3 fixed_ticker

$ws.Cells.Item(4, 4).Value = 55.21036980886064   # D4 open_price
$ws.Cells.Item(4, 5).Value = 60.91567230224609   # E4 close_price
$ws.Cells.Item(4, 6).Value = 61.58572529415808   # F4 high_price
$ws.Cells.Item(4, 7).Value = 53.74216822282875   # G4 low_price
$ws.Cells.Item(4, 8).Value = 61741703       # H4 shares_outstanding
$ws.Cells.Item(4, 9).Value = "NICE"        # I4 fixed_ticker

$ws.Cells.Item(5, 4).Value = 55.36186276160937   # D5 open_price
$ws.Cells.Item(5, 5).Value = 59.79752731323242   # E5 close_price
$ws.Cells.Item(5, 6).Value = 59.8469252797537   # F5 high_price
$ws.Cells.Item(5, 7).Value = 54.43323525792403   # G5 low_price
$ws.Cells.Item(5, 8).Value = 61741703       # H5 shares_outstanding
$ws.Cells.Item(5, 9).Value = "NICE"        # I5 fixed_ticker

$ws.Cells.Item(6, 4).Value = 63.89083375722178   # D6 open_price
$ws.Cells.Item(6, 5).Value = 63.24697113037109   # E6 close_price
$ws.Cells.Item(6, 6).Value = 64.86158153333595   # F6 high_price
$ws.Cells.Item(6, 7).Value = 60.86964004846692   # G6 low_price
$ws.Cells.Item(6, 8).Value = 61741703       # H6 shares_outstanding
$ws.Cells.Item(6, 9).Value = "NICE"        # I6 fixed_ticker

$ws.Cells.Item(7, 4).Value = 63.71711157022724   # D7 open_price
$ws.Cells.Item(7, 5).Value = 68.16549682617188   # E7 close_price
$ws.Cells.Item(7, 6).Value = 68.9697777406991   # F7 high_price
$ws.Cells.Item(7, 7).Value = 62.53551148560813   # G7 low_price
$ws.Cells.Item(7, 8).Value = 61741703       # H7 shares_outstanding
$ws.Cells.Item(7, 9).Value = "NICE"        # I7 fixed_ticker

$ws.Cells.Item(8, 4).Value = 66.75462499279209   # D8 open_price
$ws.Cells.Item(8, 5).Value = 66.15744781494141   # E8 close_price
$ws.Cells.Item(8, 6).Value = 67.35180217064277   # F8 high_price
$ws.Cells.Item(8, 7).Value = 65.6100309723584   # G8 low_price
$ws.Cells.Item(8, 8).Value = 61741703       # H8 shares_outstanding
$ws.Cells.Item(8, 9).Value = "NICE"        # I8 fixed_ticker

$ws.Cells.Item(9, 4).Value = 68.87074091629765   # D9 open_price
$ws.Cells.Item(9, 5).Value = 70.01808929443359   # E9 close_price
$ws.Cells.Item(9, 6).Value = 70.32737166326073   # F9 high_price
$ws.Cells.Item(9, 7).Value = 65.43867484237246   # G9 low_price
$ws.Cells.Item(9, 8).Value = 61741703       # H9 shares_outstanding
$ws.Cells.Item(9, 9).Value = "NICE"        # I9 fixed_ticker

$ws.Cells.Item(10, 4).Value = 68.59999847412109   # D10 open_price
$ws.Cells.Item(10, 5).Value = 67.44999694824219   # E10 close_price
$ws.Cells.Item(10, 6).Value = 68.69999694824219   # F10 high_price
$ws.Cells.Item(10, 7).Value = 66.56999969482422   # G10 low_price
$ws.Cells.Item(10, 8).Value = 61741703       # H10 shares_outstanding
$ws.Cells.Item(10, 9).Value = "NICE"        # I10 fixed_ticker

$ws.Cells.Item(11, 4).Value = 79.11000061035156   # D11 open_price
$ws.Cells.Item(11, 5).Value = 74.66000366210938   # E11 close_price
$ws.Cells.Item(11, 6).Value = 80.01000213623047   # F11 high_price
$ws.Cells.Item(11, 7).Value = 74.11000061035156   # G11 low_price
$ws.Cells.Item(11, 8).Value = 61741703       # H11 shares_outstanding
$ws.Cells.Item(11, 9).Value = "NICE"        # I11 fixed_ticker

$ws.Cells.Item(12, 4).Value = 81.48000335693359   # D12 open_price
$ws.Cells.Item(12, 5).Value = 83.30999755859375   # E12 close_price
$ws.Cells.Item(12, 6).Value = 83.77999877929688   # F12 high_price
$ws.Cells.Item(12, 7).Value = 78.48999786376953   # G12 low_price
$ws.Cells.Item(12, 8).Value = 61741703       # H12 shares_outstanding
$ws.Cells.Item(12, 9).Value = "NICE"        # I12 fixed_ticker

$ws.Cells.Item(13, 4).Value = 92.9499969482422   # D13 open_price
$ws.Cells.Item(13, 5).Value = 91.09999847412109   # E13 close_price
$ws.Cells.Item(13, 6).Value = 95.65000152587891   # F13 high_price
$ws.Cells.Item(13, 7).Value = 90.36000061035156   # G13 low_price
$ws.Cells.Item(13, 8).Value = 61741703       # H13 shares_outstanding
$ws.Cells.Item(13, 9).Value = "NICE"        # I13 fixed_ticker

$ws.Cells.Item(14, 4).Value = 93.59999847412109   # D14 open_price
$ws.Cells.Item(14, 5).Value = 95.16999816894533   # E14 close_price
$ws.Cells.Item(14, 6).Value = 98.58999633789062   # F14 high_price
$ws.Cells.Item(14, 7).Value = 88.73999786376953   # G14 low_price
$ws.Cells.Item(14, 8).Value = 61741703       # H14 shares_outstanding
$ws.Cells.Item(14, 9).Value = "NICE"        # I14 fixed_ticker

$ws.Cells.Item(15, 4).Value = 102.6100006103516   # D15 open_price
$ws.Cells.Item(15, 5).Value = 109.4000015258789   # E15 close_price
$ws.Cells.Item(15, 6).Value = 116.0500030517578   # F15 high_price
$ws.Cells.Item(15, 7).Value = 102.4100036621094   # G15 low_price
$ws.Cells.Item(15, 8).Value = 61741703       # H15 shares_outstanding
$ws.Cells.Item(15, 9).Value = "NICE"        # I15 fixed_ticker

$ws.Cells.Item(16, 4).Value = 115.3300018310547   # D16 open_price
$ws.Cells.Item(16, 5).Value = 105.9400024414062   # E16 close_price
$ws.Cells.Item(16, 6).Value = 115.3600006103516   # F16 high_price
$ws.Cells.Item(16, 7).Value = 100.5400009155273   # G16 low_price
$ws.Cells.Item(16, 8).Value = 61741703       # H16 shares_outstanding
$ws.Cells.Item(16, 9).Value = "NICE"        # I16 fixed_ticker

$ws.Cells.Item(17, 4).Value = 106.7799987792969   # D17 open_price
$ws.Cells.Item(17, 5).Value = 109.9499969482422   # E17 close_price
$ws.Cells.Item(17, 6).Value = 110.0400009155273   # F17 high_price
$ws.Cells.Item(17, 7).Value = 102.6699981689453   # G17 low_price
$ws.Cells.Item(17, 8).Value = 61741703       # H17 shares_outstanding
$ws.Cells.Item(17, 9).Value = "NICE"        # I17 fixed_ticker

$ws.Cells.Item(18, 4).Value = 124.8899993896484   # D18 open_price
$ws.Cells.Item(18, 5).Value = 137.8600006103516   # E18 close_price
$ws.Cells.Item(18, 6).Value = 140.2899932861328   # F18 high_price
$ws.Cells.Item(18, 7).Value = 121.3099975585938   # G18 low_price
$ws.Cells.Item(18, 8).Value = 61741703       # H18 shares_outstanding
$ws.Cells.Item(18, 9).Value = "NICE"        # I18 fixed_ticker

$ws.Cells.Item(19, 4).Value = 138.6300048828125   # D19 open_price
$ws.Cells.Item(19, 5).Value = 152.7200012207031   # E19 close_price
$ws.Cells.Item(19, 6).Value = 153.9900054931641   # F19 high_price
$ws.Cells.Item(19, 7).Value = 137.5599975585938   # G19 low_price
$ws.Cells.Item(19, 8).Value = 61741703       # H19 shares_outstanding
$ws.Cells.Item(19, 9).Value = "NICE"        # I19 fixed_ticker

$ws.Cells.Item(20, 4).Value = 143.4900054931641   # D20 open_price
$ws.Cells.Item(20, 5).Value = 157.7899932861328   # E20 close_price
$ws.Cells.Item(20, 6).Value = 159.3000030517578   # F20 high_price
$ws.Cells.Item(20, 7).Value = 139.0599975585938   # G20 low_price
$ws.Cells.Item(20, 8).Value = 61741703       # H20 shares_outstanding
$ws.Cells.Item(20, 9).Value = "NICE"        # I20 fixed_ticker

$ws.Cells.Item(21, 4).Value = 156.7200012207031   # D21 open_price
$ws.Cells.Item(21, 5).Value = 172.3000030517578   # E21 close_price
$ws.Cells.Item(21, 6).Value = 182.0599975585937   # F21 high_price
$ws.Cells.Item(21, 7).Value = 156.0099945068359   # G21 low_price
$ws.Cells.Item(21, 8).Value = 61741703       # H21 shares_outstanding
$ws.Cells.Item(21, 9).Value = "NICE"        # I21 fixed_ticker

$ws.Cells.Item(22, 4).Value = 141.9799957275391   # D22 open_price
$ws.Cells.Item(22, 5).Value = 164.3000030517578   # E22 close_price
$ws.Cells.Item(22, 6).Value = 170.7899932861328   # F22 high_price
$ws.Cells.Item(22, 7).Value = 139.8800048828125   # G22 low_price
$ws.Cells.Item(22, 8).Value = 61741703       # H22 shares_outstanding
$ws.Cells.Item(22, 9).Value = "NICE"        # I22 fixed_ticker

$ws.Cells.Item(23, 4).Value = 187.8500061035156   # D23 open_price
$ws.Cells.Item(23, 5).Value = 205.2400054931641   # E23 close_price
$ws.Cells.Item(23, 6).Value = 206.5   # F23 high_price
$ws.Cells.Item(23, 7).Value = 186.1600036621093   # G23 low_price
$ws.Cells.Item(23, 8).Value = 61741703       # H23 shares_outstanding
$ws.Cells.Item(23, 9).Value = "NICE"        # I23 fixed_ticker

$ws.Cells.Item(24, 4).Value = 228.2700042724609   # D24 open_price
$ws.Cells.Item(24, 5).Value = 228.259994506836   # E24 close_price
$ws.Cells.Item(24, 6).Value = 240.5800018310547   # F24 high_price
$ws.Cells.Item(24, 7).Value = 223.3000030517578   # G24 low_price
$ws.Cells.Item(24, 8).Value = 61741703       # H24 shares_outstanding
$ws.Cells.Item(24, 9).Value = "NICE"        # I24 fixed_ticker

$ws.Cells.Item(25, 4).Value = 281.510009765625   # D25 open_price
$ws.Cells.Item(25, 5).Value = 261.2799987792969   # E25 close_price
$ws.Cells.Item(25, 6).Value = 282.2900085449219   # F25 high_price
$ws.Cells.Item(25, 7).Value = 249.7200012207031   # G25 low_price
$ws.Cells.Item(25, 8).Value = 61741703       # H25 shares_outstanding
$ws.Cells.Item(25, 9).Value = "NICE"        # I25 fixed_ticker

$ws.Cells.Item(26, 4).Value = 218.8600006103516   # D26 open_price
$ws.Cells.Item(26, 5).Value = 241.229995727539   # E26 close_price
$ws.Cells.Item(26, 6).Value = 244.4600067138672   # F26 high_price
$ws.Cells.Item(26, 7).Value = 218.8600006103516   # G26 low_price
$ws.Cells.Item(26, 8).Value = 61741703       # H26 shares_outstanding
$ws.Cells.Item(26, 9).Value = "NICE"        # I26 fixed_ticker

$ws.Cells.Item(27, 4).Value = 247.4700012207031   # D27 open_price
$ws.Cells.Item(27, 5).Value = 278.6499938964844   # E27 close_price
$ws.Cells.Item(27, 6).Value = 287.8500061035156   # F27 high_price
$ws.Cells.Item(27, 7).Value = 246.9700012207031   # G27 low_price
$ws.Cells.Item(27, 8).Value = 61741703       # H27 shares_outstanding
$ws.Cells.Item(27, 9).Value = "NICE"        # I27 fixed_ticker

$ws.Cells.Item(28, 4).Value = 284.8200073242188   # D28 open_price
$ws.Cells.Item(28, 5).Value = 283.0199890136719   # E28 close_price
$ws.Cells.Item(28, 6).Value = 285.6300048828125   # F28 high_price
$ws.Cells.Item(28, 7).Value = 262.0899963378906   # G28 low_price
$ws.Cells.Item(28, 8).Value = 61741703       # H28 shares_outstanding
$ws.Cells.Item(28, 9).Value = "NICE"        # I28 fixed_ticker

$ws.Cells.Item(29, 4).Value = 303.75   # D29 open_price
$ws.Cells.Item(29, 5).Value = 256.0599975585937   # E29 close_price
$ws.Cells.Item(29, 6).Value = 306.1199951171875   # F29 high_price
$ws.Cells.Item(29, 7).Value = 232.9499969482422   # G29 low_price
$ws.Cells.Item(29, 8).Value = 61741703       # H29 shares_outstanding
$ws.Cells.Item(29, 9).Value = "NICE"        # I29 fixed_ticker

$ws.Cells.Item(30, 4).Value = 229.9900054931641   # D30 open_price
$ws.Cells.Item(30, 5).Value = 206.4100036621093   # E30 close_price
$ws.Cells.Item(30, 6).Value = 234.4900054931641   # F30 high_price
$ws.Cells.Item(30, 7).Value = 200.9199981689453   # G30 low_price
$ws.Cells.Item(30, 8).Value = 61741703       # H30 shares_outstanding
$ws.Cells.Item(30, 9).Value = "NICE"        # I30 fixed_ticker

$ws.Cells.Item(31, 4).Value = 192.1799926757812   # D31 open_price
$ws.Cells.Item(31, 5).Value = 214.0200042724609   # E31 close_price
$ws.Cells.Item(31, 6).Value = 215.259994506836   # F31 high_price
$ws.Cells.Item(31, 7).Value = 190.4400024414062   # G31 low_price
$ws.Cells.Item(31, 8).Value = 61741703       # H31 shares_outstanding
$ws.Cells.Item(31, 9).Value = "NICE"        # I31 fixed_ticker

$ws.Cells.Item(32, 4).Value = 187.8600006103516   # D32 open_price
$ws.Cells.Item(32, 5).Value = 189.8899993896484   # E32 close_price
$ws.Cells.Item(32, 6).Value = 204.2799987792969   # F32 high_price
$ws.Cells.Item(32, 7).Value = 178.2799987792969   # G32 low_price
$ws.Cells.Item(32, 8).Value = 61741703       # H32 shares_outstanding
$ws.Cells.Item(32, 9).Value = "NICE"        # I32 fixed_ticker

$ws.Cells.Item(33, 4).Value = 192.3000030517578   # D33 open_price
$ws.Cells.Item(33, 5).Value = 207.4299926757812   # E33 close_price
$ws.Cells.Item(33, 6).Value = 214.7700042724609   # F33 high_price
$ws.Cells.Item(33, 7).Value = 182.8899993896484   # G33 low_price
$ws.Cells.Item(33, 8).Value = 61741703       # H33 shares_outstanding
$ws.Cells.Item(33, 9).Value = "NICE"        # I33 fixed_ticker

$ws.Cells.Item(34, 4).Value = 226.5399932861328   # D34 open_price
$ws.Cells.Item(34, 5).Value = 204.009994506836   # E34 close_price
$ws.Cells.Item(34, 6).Value = 229.8699951171875   # F34 high_price
$ws.Cells.Item(34, 7).Value = 198.4700012207031   # G34 low_price
$ws.Cells.Item(34, 8).Value = 61741703       # H34 shares_outstanding
$ws.Cells.Item(34, 9).Value = "NICE"        # I34 fixed_ticker

$ws.Cells.Item(35, 4).Value = 203.3500061035156   # D35 open_price
$ws.Cells.Item(35, 5).Value = 217.8500061035156   # E35 close_price
$ws.Cells.Item(35, 6).Value = 225.3699951171875   # F35 high_price
$ws.Cells.Item(35, 7).Value = 195.9100036621093   # G35 low_price
$ws.Cells.Item(35, 8).Value = 61741703       # H35 shares_outstanding
$ws.Cells.Item(35, 9).Value = "NICE"        # I35 fixed_ticker

$ws.Cells.Item(36, 4).Value = 170.5500030517578   # D36 open_price
$ws.Cells.Item(36, 5).Value = 154.3500061035156   # E36 close_price
$ws.Cells.Item(36, 6).Value = 174.1799926757812   # F36 high_price
$ws.Cells.Item(36, 7).Value = 150.1799926757812   # G36 low_price
$ws.Cells.Item(36, 8).Value = 61741703       # H36 shares_outstanding
$ws.Cells.Item(36, 9).Value = "NICE"        # I36 fixed_ticker

$ws.Cells.Item(37, 4).Value = 198.7700042724609   # D37 open_price
$ws.Cells.Item(37, 5).Value = 208.1000061035156   # E37 close_price
$ws.Cells.Item(37, 6).Value = 220.4400024414062   # F37 high_price
$ws.Cells.Item(37, 7).Value = 191.3099975585937   # G37 low_price
$ws.Cells.Item(37, 8).Value = 61741703       # H37 shares_outstanding
$ws.Cells.Item(37, 9).Value = "NICE"        # I37 fixed_ticker

$ws.Cells.Item(38, 4).Value = 260.0199890136719   # D38 open_price
$ws.Cells.Item(38, 5).Value = 223.509994506836   # E38 close_price
$ws.Cells.Item(38, 6).Value = 260.6199951171875   # F38 high_price
$ws.Cells.Item(38, 7).Value = 222.1399993896484   # G38 low_price
$ws.Cells.Item(38, 8).Value = 61741703       # H38 shares_outstanding
$ws.Cells.Item(38, 9).Value = "NICE"        # I38 fixed_ticker

$ws.Cells.Item(39, 4).Value = 172.6699981689453   # D39 open_price
$ws.Cells.Item(39, 5).Value = 181   # E39 close_price
$ws.Cells.Item(39, 6).Value = 185.4100036621093   # F39 high_price
$ws.Cells.Item(39, 7).Value = 167.5200042724609   # G39 low_price
$ws.Cells.Item(39, 8).Value = 61741703       # H39 shares_outstanding
$ws.Cells.Item(39, 9).Value = "NICE"        # I39 fixed_ticker

$ws.Cells.Item(40, 4).Value = 172.4700012207031   # D40 open_price
$ws.Cells.Item(40, 5).Value = 173.6999969482422   # E40 close_price
$ws.Cells.Item(40, 6).Value = 182.9700012207031   # F40 high_price
$ws.Cells.Item(40, 7).Value = 162.5299987792969   # G40 low_price
$ws.Cells.Item(40, 8).Value = 61741703       # H40 shares_outstanding
$ws.Cells.Item(40, 9).Value = "NICE"        # I40 fixed_ticker

$ws.Cells.Item(41, 4).Value = 171.6300048828125   # D41 open_price
$ws.Cells.Item(41, 5).Value = 166.1399993896484   # E41 close_price
$ws.Cells.Item(41, 6).Value = 173   # F41 high_price
$ws.Cells.Item(41, 7).Value = 158.5099945068359   # G41 low_price
$ws.Cells.Item(41, 8).Value = 61741703       # H41 shares_outstanding
$ws.Cells.Item(41, 9).Value = "NICE"        # I41 fixed_ticker

$ws.Cells.Item(42, 4).Value = 155.2899932861328   # D42 open_price
$ws.Cells.Item(42, 5).Value = 155.8500061035156   # E42 close_price
$ws.Cells.Item(42, 6).Value = 158.4700012207031   # F42 high_price
$ws.Cells.Item(42, 7).Value = 138.7899932861328   # G42 low_price
$ws.Cells.Item(42, 8).Value = 61741703       # H42 shares_outstanding
$ws.Cells.Item(42, 9).Value = "NICE"        # I42 fixed_ticker

$ws.Cells.Item(43, 4).Value = 169.1499938964844   # D43 open_price
$ws.Cells.Item(43, 5).Value = 156.0500030517578   # E43 close_price
$ws.Cells.Item(43, 6).Value = 175   # F43 high_price
$ws.Cells.Item(43, 7).Value = 150.0800018310547   # G43 low_price
$ws.Cells.Item(43, 8).Value = 61741703       # H43 shares_outstanding
$ws.Cells.Item(43, 9).Value = "NICE"        # I43 fixed_ticker
